# Restore revision #7341d9fba0d6f2f520139e28f4ba3083eeb6451c.TEST
# The "Rules" sheet holds hour-range rules for a greeting function.
# Rule R30 ("Good Evening") had its "From" (Integer min) threshold
# stored in C10. This revision changes that threshold from 18 back to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
